$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 13.125
$ws.Range("C3").Value = -11.817
$ws.Range("A4").Value = -21.573
$ws.Range("B4").Value = 6.962999999999999
$ws.Range("C4").Value = -12.585
$ws.Range("B5").Value = 6.578
$ws.Range("A6").Value = -21.351
$ws.Range("A7").Value = -21.361
$ws.Range("B8").Value = 6.153
$ws.Range("C9").Value = -11.775
$ws.Range("C11").Value = -12.642
$ws.Range("C14").Value = -11.607
$ws.Range("A16").Value = -21.361
$ws.Range("B16").Value = 5.896
$ws.Range("C18").Value = -12.412
$ws.Range("A20").Value = -22.2
$ws.Range("E20").Value = 12.951
$ws.Range("B22").Value = 6.693000000000001
$ws.Range("C25").Value = -12.642
